$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns to snake_case names for use in R
$ws.Range("A1").Value = "patient"
$ws.Range("B1").Value = "mood_pre"
$ws.Range("C1").Value = "mood_post"

# Reset the active selection back to A1 (clears the saved B2 selection state)
$ws.Range("A1").Select()
